$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = 117591.8247312488
$ws.Range("E3").Value = -0.03063365246675705
$ws.Range("F3").Value = 0.218308688172155
$ws.Range("G3").Value = -2.349930543709918
$ws.Range("H3").Value = 31.76095076488103

$ws.Range("D4").Value = 118239.923232096
$ws.Range("E4").Value = -0.03039393502373985
$ws.Range("F4").Value = 0.2232792613453317
$ws.Range("G4").Value = -1.737613789996416
$ws.Range("H4").Value = 21.14839510600624

$ws.Range("D5").Value = 118966.0236033464
$ws.Range("E5").Value = -0.03471270622512351
$ws.Range("F5").Value = 0.2317978979356952
$ws.Range("G5").Value = -1.202008007481708
$ws.Range("H5").Value = 12.32410070557436

$ws.Range("D6").Value = 119737.9283833025
$ws.Range("E6").Value = -0.03325136278406547
$ws.Range("F6").Value = 0.2203887710611357
$ws.Range("G6").Value = -0.5977070996180361
$ws.Range("H6").Value = 7.668406080852271

$ws.Range("D7").Value = 120283.1691578465
$ws.Range("E7").Value = -0.04132931447787107
$ws.Range("F7").Value = 0.2594842032796034
$ws.Range("G7").Value = -1.752565880427223
$ws.Range("H7").Value = 16.69005240226756

$ws.Range("D8").Value = 122675.1844972577
$ws.Range("E8").Value = -0.05089586534709097
$ws.Range("F8").Value = 0.2283332993542903
$ws.Range("G8").Value = -1.041006549741057
$ws.Range("H8").Value = 8.050611020791827

$ws.Range("D10").Value = 124827.8243651185
$ws.Range("E10").Value = -0.1264793875796638
$ws.Range("F10").Value = 0.4559626492714848
$ws.Range("G10").Value = -1.956310446760939
$ws.Range("H10").Value = 9.784302688364713

$ws.Range("D11").Value = 126897.7148185333
$ws.Range("E11").Value = -0.1385790421169318
$ws.Range("F11").Value = 0.4508392737052024
$ws.Range("G11").Value = -1.670682921823281
$ws.Range("H11").Value = 7.793436157897712

$ws.Range("D15").Value = 115963.2299981778
$ws.Range("E15").Value = -0.09635326435678591
$ws.Range("F15").Value = 0.144233387973362
$ws.Range("G15").Value = -0.8875382030323288
$ws.Range("H15").Value = 4.027310534457601

$ws.Range("D17").Value = 116724.1843808662
$ws.Range("E17").Value = -0.07372048939007186
$ws.Range("F17").Value = 0.205799471962948
$ws.Range("G17").Value = -0.1647116387313506
$ws.Range("H17").Value = 6.099566578310869

$ws.Range("D19").Value = 117657.9442087174
$ws.Range("E19").Value = -0.01609327286070889
$ws.Range("F19").Value = 0.1847166487167634
$ws.Range("G19").Value = -0.03157244337971703
$ws.Range("H19").Value = 4.689487087389645
